$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new rows (4 and 5) replicating the pattern of row 3, with new
# timestamps in column A. This mirrors additional price-check records
# captured (e.g. by the new HtmlHelper class) for GILD.

$ws.Range("A4").Value = 42601.882060185184
$ws.Range("B4").Value = "Gilead Sciences, Inc."
$ws.Range("C4").Value = "GILD"
$ws.Range("D4").Value = 80.91
$ws.Range("E4").Value = 80.7
$ws.Range("F4").Value = -0.14000000000000001

$ws.Range("A5").Value = 42601.883564814816
$ws.Range("B5").Value = "Gilead Sciences, Inc."
$ws.Range("C5").Value = "GILD"
$ws.Range("D5").Value = 80.91
$ws.Range("E5").Value = 80.7
$ws.Range("F5").Value = -0.14000000000000001

# Match the date/time number format used by the existing date column (A2:A3)
$ws.Range("A4:A5").NumberFormat = "m/d/yy h:mm"
